# Applies the "Saldo" export update:
#   - account 004749680/OTAVIANO (row 4) is replaced by 004001621/DANIELA
#     with a new balance of 100484.68, immediately followed by a brand new
#     row for 003301389/EDMUNDO (70000)
#   - the old 004001621/DANIELA row (balance 484.68) further down the sheet
#     is removed (its balance was folded into the new row above)
#   - the 005055865/G3C row is removed entirely
#   - six new accounts are inserted right after 004224011/THOMAS:
#     EULER, VERA, DAIANNE, PEDRO, MARIA, GUILHERME
#   - one new account (005348975/JULIA) is inserted right after the
#     005266369/EG row
#
# All row numbers below refer to the ORIGINAL (before-edit) layout and the
# edits are applied from the bottom of the sheet upward so that an
# insert/delete never invalidates a row number used by a later step.
#
# Account numbers are zero-padded numeric strings (e.g. "004001621"), so
# each "Conta" cell is explicitly forced to text format before its value
# is assigned - otherwise Excel would coerce it to a number and drop the
# leading zeros.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-AccountRow($row, $conta, $nome, $saldo) {
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $conta
    $ws.Cells.Item($row, 2).Value = $nome
    $ws.Cells.Item($row, 3).Value = $saldo
}

# --- 1. insert JULIA right after the EG row (original row 84) ----------
$ws.Rows.Item(85).Insert()
Set-AccountRow 85 "005348975" "JULIA" 400

# --- 2. delete the old DANIELA row (004001621 / 484.68, original row 71) -
$ws.Rows.Item(71).Delete()

# --- 3. insert six new accounts right after THOMAS (original row 7) ----
$ws.Rows.Item(8).Insert()
Set-AccountRow 8 "004399832" "EULER" 16614.4

$ws.Rows.Item(9).Insert()
Set-AccountRow 9 "005105970" "VERA" 15000

$ws.Rows.Item(10).Insert()
Set-AccountRow 10 "004473942" "DAIANNE" 6248.86

$ws.Rows.Item(11).Insert()
Set-AccountRow 11 "004460491" "PEDRO" 4988.96

$ws.Rows.Item(12).Insert()
Set-AccountRow 12 "004870019" "MARIA" 3805.81

$ws.Rows.Item(13).Insert()
Set-AccountRow 13 "004574428" "GUILHERME" 2285.83

# --- 4. delete the G3C row (original row 6) -----------------------------
$ws.Rows.Item(6).Delete()

# --- 5. turn the OTAVIANO row (original row 4) into DANIELA ------------
Set-AccountRow 4 "004001621" "DANIELA" 100484.68

# --- 6. insert the brand new EDMUNDO row right after it -----------------
$ws.Rows.Item(5).Insert()
Set-AccountRow 5 "003301389" "EDMUNDO" 70000
